# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Only column G ("K") values change; everything else in the sheet (dates,
# TB/PC/dS0/dSF/IP/I0/IF columns, headers, styles) stays as-is. Write the
# recalculated K values for rows 2-20 (rows 17 and 21 already hold their
# correct value and are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 5
    4  = 1
    5  = 4
    6  = 4
    7  = 4
    8  = 5
    9  = 5
    10 = 6
    11 = 7
    12 = 2
    13 = 4
    14 = 4
    15 = 0
    16 = 3
    18 = 4
    19 = 5
    20 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
